# Refresh the cryptos price/volume snapshot (GitHub Actions style update).
# Price (column D) values are entered with a leading apostrophe so Excel
# keeps them as literal text (e.g. "1.005") instead of re-parsing them as
# numbers/dates, matching the original inline-string cell contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.597.85'
$ws.Range("E2").Value = '  +2.00%  '
$ws.Range("D3").Value = '''1.868.65'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("D4").Value = '''1.006'
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").Value = '''325.89'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = '''0.4650'
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").Value = '''0.3880'
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("D9").Value = '''0.07860'
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '''0.9726'
$ws.Range("E10").Value = '  +1.28%  '
$ws.Range("D11").Value = '''21.92'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '''1.853.92'
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").Value = '''6.979'
$ws.Range("E13").Value = '  +1.24%  '
$ws.Range("D14").Value = '''5.692'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").Value = '''0.06981'
$ws.Range("E15").Value = '  +3.64%  '
$ws.Range("D16").Value = '''88.00'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '''1.005'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '''0.00001003'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").Value = '''16.80'
$ws.Range("E19").Value = '  +1.03%  '
$ws.Range("D20").Value = '''1.004'
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").Value = '''28.618.41'
$ws.Range("E21").Value = '  +2.03%  '
$ws.Range("D22").Value = '''5.283'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").Value = '''10.99'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = '''2.116'
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("D25").Value = '''2.144.56'
$ws.Range("E25").Value = '  +2.88%  '
$ws.Range("D26").Value = '''152.67'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").Value = '''19.20'
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").Value = '''5.781'
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("D29").Value = '''1.985'
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("D30").Value = '''119.37'
$ws.Range("E30").Value = '  +1.82%  '
$ws.Range("D31").Value = '''0.09363'
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").Value = '''0.9193'
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("D33").Value = '''5.259'
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("D34").Value = '''1.336'
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").Value = '''3.335'
$ws.Range("E35").Value = '  +0.54%  '
$ws.Range("D36").Value = '''0.05795'
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("D37").Value = '''0.02095'
$ws.Range("E37").Value = '  -2.21%  '
$ws.Range("D38").Value = '''1.148'
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").Value = '''7.752'
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("D40").Value = '''0.5619'
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("D41").Value = '''0.1785'
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("E42").Value = '  -1.45%  '
$ws.Range("D43").Value = '''0.07204'
$ws.Range("E43").Value = '  +2.59%  '
$ws.Range("D44").Value = '''11.76'
$ws.Range("E44").Value = '  +1.55%  '
$ws.Range("D45").Value = '''0.5308'
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").Value = '''1.145'
$ws.Range("E46").Value = '  -4.97%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''2.084'
$ws.Range("E47").Value = '  -3.41%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''1.822'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("D49").Value = '''113.14'
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("D50").Value = '''2.407'
$ws.Range("E50").Value = '  +3.88%  '
$ws.Range("D51").Value = '''1.003'
$ws.Range("E51").Value = '  +0.28%  '
